$wb = $excel.ActiveWorkbook

# --- Data edit: rename "质控组" -> "北京组" group label ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2:A5").Value = "北京组"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2").Value = "北京组"

# --- View state: Sheet1 becomes the active sheet/tab, with A5 selected ---
[void]$ws1.Activate()
[void]$ws1.Range("A5").Select()

# Sheet2 keeps a lingering selection at A2 (no longer the active tab)
[void]$ws2.Range("A2").Select()

# Re-activate Sheet1 so it is the tab shown when the workbook is saved/reopened
[void]$ws1.Activate()
